$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '{''chika'', ''$'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D3").Value = '{''chika'', ''$'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D4").Value = '{''chika'', ''$'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D5").Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool'', ''naur''}'
$ws.Range("D6").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D8").Value = '{''chika'', ''$'', ''naur'', ''shimenet'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D10").Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool''}'
$ws.Range("D11").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D18").Value = '{'')'', ''('', '','', ''['', ''='', '';''}'
$ws.Range("D19").Value = '{'')'', ''('', '','', ''['', ''='', '';''}'
$ws.Range("D20").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D22").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D26").Value = '{'')'', '','', ''('', ''='', '';''}'
$ws.Range("D28").Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''='', ''+'', ''+='', ''step'', ''**='', ''!='', ''/='', ''%='', '';'', ''<='', ''//='', ''*='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', ''-='', '')'', ''<'', ''>'', ''}''}'
$ws.Range("D30").Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''='', ''+'', ''+='', ''step'', ''**='', ''!='', ''/='', ''%='', '';'', ''<='', ''//='', ''*='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', ''-='', '')'', ''<'', ''>'', ''}''}'
$ws.Range("D32").Value = '{''<='', ''//='', ''*='', ''||'', ''to'', ''>='', ''-'', ''*'', ''**'', ''=='', '']'', ''/'', ''='', ''&&'', ''+'', ''+='', ''step'', ''//'', ''%'', '','', '':'', ''-='', ''**='', ''!='', ''/='', '')'', ''<'', ''%='', ''>'', ''}'', '';''}'
$ws.Range("D33").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D40").Value = '{''<='', ''||'', ''to'', ''>='', ''('', ''['', ''-'', ''*'', ''**'', ''=='', '']'', ''/'', ''&&'', ''+'', ''step'', ''//'', ''%'', '','', '':'', ''!='', '')'', ''<'', ''>'', ''}'', '';''}'
$ws.Range("D41").Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range("D43").Value = '{''||'', ''to'', ''['', ''-'', ''**'', ''=='', ''='', ''+'', ''+='', ''step'', ''**='', ''!='', ''/='', ''%='', '';'', ''<='', ''//='', ''*='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', ''-='', '')'', ''<'', ''>'', ''}''}'
$ws.Range("D46").Value = '{'','', '';''}'
$ws.Range("D50").Value = '{'','', '';''}'
$ws.Range("D51").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D55").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("B60").Value = '<func-def> -> <return-type> id <array-dec> ( <parameters> ) { <func-body> } <func-def>'
$ws.Range("D60").Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool'', ''shimenet''}'
$ws.Range("D62").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D64").Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool'', ''naur''}'
$ws.Range("D65").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D66").Value = '{''chika'', ''anda'', ''andamhie'', ''eklabool'', ''naur''}'
$ws.Range("D67").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D68").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D72").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D74").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D75").Value = '{''push'', ''keri'', ''naur'', ''++'', ''--'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D76").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D78").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D80").Value = '{''<='', ''//'', ''%'', ''>='', ''||'', ''-'', ''!='', ''*'', ''**'', ''<'', ''=='', ''>'', ''/'', ''&&'', ''+''}'
$ws.Range("D81").Value = '{''step'', '','', '':'', ''to'', '')'', '']'', ''}'', '';''}'
$ws.Range("D82").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D85").Value = '{''eme'', ''len'', ''('', ''++'', ''--'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D88").Value = '{''eme'', ''chika_literal'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D93").Value = '{''<='', ''||'', ''to'', ''>='', ''('', ''['', ''-'', ''*'', ''**'', ''=='', '']'', ''/'', ''&&'', ''+'', ''step'', ''//'', ''%'', '','', '':'', ''!='', '')'', ''<'', ''>'', ''}'', '';''}'
$ws.Range("D95").Value = '{''||'', ''to'', ''-'', ''**'', ''=='', ''+'', ''step'', ''!='', '';'', ''<='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', '')'', ''<'', ''>'', ''}''}'
$ws.Range("D98").Value = '{''eme'', ''chika_literal'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D115").Value = '{''chika'', ''anda'', ''andamhie'', ''id'', ''eklabool''}'
$ws.Range("D120").Value = '{''forda'', ''keri''}'
$ws.Range("D122").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''ditech'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D123").Value = '{''++'', ''id'', ''--''}'
$ws.Range("D124").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D125").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D127").Value = '{''chika'', ''anda'', ''andamhie'', ''id'', ''eklabool''}'
$ws.Range("D131").Value = '{''forda'', ''keri''}'
$ws.Range("D134").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''ditech'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D135").Value = '{''++'', ''id'', ''--''}'
$ws.Range("D136").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''ditech'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D137").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D138").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D139").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D141").Value = '{''||'', ''to'', ''['', ''-'', ''**'', ''=='', ''='', ''+'', ''+='', ''step'', ''**='', ''!='', ''/='', ''%='', '';'', ''<='', ''//='', ''*='', ''>='', ''*'', '']'', ''/'', ''&&'', ''//'', ''%'', '','', '':'', ''-='', '')'', ''<'', ''>'', ''}''}'
$ws.Range("D151").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D153").Value = '{''chika'', ''anda'', ''andamhie'', ''id'', ''eklabool''}'
$ws.Range("D154").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D156").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D160").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D162").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D163").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D166").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D169").Value = '{''push'', ''ganern'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D171").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D175").Value = '{''andamhie'', ''eklabool'', ''chika'', ''anda''}'
$ws.Range("D177").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D178").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D181").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D183").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D184").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D189").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D193").Value = '{''}'', ''ditech''}'
$ws.Range("D194").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
$ws.Range("D195").Value = '{''amaccana'', ''push'', ''gogogo'', ''keri'', ''naur'', ''++'', ''--'', ''betsung'', ''adele'', ''adelete'', ''eklabool'', ''forda'', ''chika'', ''ditech'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''andamhie'', ''id''}'
$ws.Range("D198").Value = '{''}'', ''ditech''}'
$ws.Range("D199").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D201").Value = '{''ditech'', ''betsung'', ''}''}'
$ws.Range("D206").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D208").Value = '{''push'', ''gogogo'', ''++'', ''--'', ''adele'', ''adelete'', ''ditech'', ''andamhie'', ''amaccana'', ''keri'', ''naur'', ''betsung'', ''eklabool'', ''forda'', ''chika'', ''serve'', ''pak'', ''versa'', ''anda'', ''}'', ''id''}'
$ws.Range("D211").Value = '{''eme'', ''len'', ''-'', ''('', ''++'', ''--'', ''!'', ''chika_literal'', ''id'', ''korik'', ''anda_literal'', ''andamhie_literal''}'
